$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 421.66666
$ws.Range("I28").Value = 291.5625
$ws.Range("J28").Value = 838
$ws.Range("K28").Value = 291.5625
$ws.Range("L28").Value = 838
$ws.Range("M28").Value = 193.4375
$ws.Range("N28").Value = -1808

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H70").Value = 1580.5625
$ws.Range("I70").Value = 1248.375
$ws.Range("J70").Value = 1912.75
$ws.Range("K70").Value = 3745.125
$ws.Range("L70").Value = 5738.25
$ws.Range("M70").Value = -3475.125
$ws.Range("N70").Value = -6278.25

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H73").Value = 1580.5625
$ws.Range("I73").Value = 1248.375
$ws.Range("J73").Value = 1912.75
$ws.Range("K73").Value = 3745.125
$ws.Range("L73").Value = 5738.25
$ws.Range("M73").Value = -2809.125
$ws.Range("N73").Value = -7610.25

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H80").Value = 4168.303
$ws.Range("I80").Value = 822.06665
$ws.Range("J80").Value = 6956.8335
$ws.Range("K80").Value = 2466.19995
$ws.Range("L80").Value = 20870.5005
$ws.Range("M80").Value = -1468.19995
$ws.Range("N80").Value = -22866.5005

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H83").Value = 4168.303
$ws.Range("I83").Value = 822.06665
$ws.Range("J83").Value = 6956.8335
$ws.Range("K83").Value = 7398.59985
$ws.Range("L83").Value = 62611.5015
$ws.Range("M83").Value = -2406.59985
$ws.Range("N83").Value = -72595.5015

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H86").Value = 56914.05
$ws.Range("I86").Value = 128624.125
$ws.Range("J86").Value = 4761.273
$ws.Range("K86").Value = 128624.125
$ws.Range("L86").Value = 4761.273
$ws.Range("M86").Value = -127501.125
$ws.Range("N86").Value = -7007.273

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H88").Value = 1540.3334
$ws.Range("J88").Value = 1576.1111
$ws.Range("L88").Value = 1576.1111
$ws.Range("N88").Value = -2388.1111

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H89").Value = 56914.05
$ws.Range("I89").Value = 128624.125
$ws.Range("J89").Value = 4761.273
$ws.Range("K89").Value = 643120.625
$ws.Range("L89").Value = 23806.365
$ws.Range("M89").Value = -637504.625
$ws.Range("N89").Value = -35038.36500000001

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H91").Value = 1540.3334
$ws.Range("J91").Value = 1576.1111
$ws.Range("L91").Value = 1576.1111
$ws.Range("N91").Value = -4384.1111

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H126").Value = 0
$ws.Range("J126").Value = 0
$ws.Range("L126").ClearContents()
$ws.Range("N126").Value = 0

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 1989.7435
$ws.Range("I132").Value = 1814.1177
$ws.Range("J132").Value = 3184
$ws.Range("K132").Value = 5442.3531
$ws.Range("L132").Value = 9552
$ws.Range("M132").Value = -2912.3531
$ws.Range("N132").Value = -14612

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 3924.2144
$ws.Range("I138").Value = 1522.8636
$ws.Range("J138").Value = 5024.8335
$ws.Range("K138").Value = 4568.5908
$ws.Range("L138").Value = 15074.5005
$ws.Range("M138").Value = 571.4092000000001
$ws.Range("N138").Value = -25354.5005

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H141").Value = 2799.8125
$ws.Range("I141").Value = 2241.5789
$ws.Range("J141").Value = 3615.6924
$ws.Range("K141").Value = 6724.736699999999
$ws.Range("L141").Value = 10847.0772
$ws.Range("M141").Value = -1544.736699999999
$ws.Range("N141").Value = -21207.0772

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1987.1852
$ws.Range("I2").Value = 2087.15
$ws.Range("K2").Value = 2087.15
$ws.Range("M2").Value = -1974.15

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 650
$ws.Range("I4").Value = 0
$ws.Range("K4").Value = 0
$ws.Range("M4").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H9").Value = 40009
$ws.Range("J9").Value = 40009
$ws.Range("L9").Value = 40009
$ws.Range("N9").Value = -40349

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H20").Value = 40009
$ws.Range("J20").Value = 40009
$ws.Range("L20").Value = 40009
$ws.Range("N20").Value = -40549

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H23").Value = 38175.09
$ws.Range("J23").Value = 28879.334
$ws.Range("L23").Value = 28879.334
$ws.Range("N23").Value = -29397.334

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H37").Value = 8000
$ws.Range("I37").Value = 8000
$ws.Range("J37").Value = 0
$ws.Range("K37").Value = 8000
$ws.Range("L37").Value = 0
$ws.Range("M37").ClearContents()
$ws.Range("N37").Value = -7727

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H44").Value = 200040800
$ws.Range("J44").Value = 200040800
$ws.Range("L44").Value = 200040800
$ws.Range("N44").Value = -200041776

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H55").Value = 59125
$ws.Range("J55").Value = 59125
$ws.Range("L55").Value = 59125
$ws.Range("N55").Value = -59755

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H116").Value = 1987.1852
$ws.Range("I116").Value = 2087.15
$ws.Range("K116").Value = 2087.15
$ws.Range("M116").Value = 206.8499999999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 6251900.5
$ws.Range("I122").Value = 1894.1333
$ws.Range("J122").Value = 25001920
$ws.Range("K122").Value = 5682.3999
$ws.Range("L122").Value = 75005760
$ws.Range("M122").Value = -3232.3999
$ws.Range("N122").Value = -75010660

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1987.1852
$ws.Range("I3").Value = 2087.15
$ws.Range("K3").Value = 2087.15
$ws.Range("M3").Value = -1973.15

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 40.333332
$ws.Range("I22").Value = 31
$ws.Range("J22").Value = 45
$ws.Range("K22").Value = 31
$ws.Range("L22").Value = 45
$ws.Range("M22").Value = 142
$ws.Range("N22").Value = -391

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H60").Value = 49995
$ws.Range("J60").Value = 49995
$ws.Range("L60").Value = 49995
$ws.Range("N60").Value = -51193

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H139").Value = 54900
$ws.Range("I139").Value = 0
$ws.Range("J139").Value = 54900
$ws.Range("K139").Value = 0
$ws.Range("L139").ClearContents()
$ws.Range("M139").Value = 54900
$ws.Range("N139").Value = -65180

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 2276095.5
$ws.Range("I58").Value = 4134822.5
$ws.Range("J58").Value = 4317.778
$ws.Range("K58").Value = 4134822.5
$ws.Range("L58").Value = 4317.778
$ws.Range("M58").Value = -4134619.5
$ws.Range("N58").Value = -4723.778

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H136").Value = 2276095.5
$ws.Range("I136").Value = 4134822.5
$ws.Range("J136").Value = 4317.778
$ws.Range("K136").Value = 12404467.5
$ws.Range("L136").Value = 12953.334
$ws.Range("M136").Value = -12401917.5
$ws.Range("N136").Value = -18053.334

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 8931305
$ws.Range("I5").Value = 553.6667
$ws.Range("J5").Value = 25006658
$ws.Range("K5").Value = 1661.0001
$ws.Range("L5").Value = 75019974
$ws.Range("M5").Value = -1549.0001
$ws.Range("N5").Value = -75020198

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H50").Value = 265.52942
$ws.Range("I50").Value = 115.28571
$ws.Range("J50").Value = 370.7
$ws.Range("K50").Value = 345.85713
$ws.Range("L50").Value = 1112.1
$ws.Range("M50").Value = 135.14287
$ws.Range("N50").Value = -2074.1

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H53").Value = 265.52942
$ws.Range("I53").Value = 115.28571
$ws.Range("J53").Value = 370.7
$ws.Range("K53").Value = 345.85713
$ws.Range("L53").Value = 1112.1
$ws.Range("M53").Value = 135.14287
$ws.Range("N53").Value = -2074.1

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H54").Value = 2979.75
$ws.Range("I54").Value = 0
$ws.Range("J54").Value = 2979.75
$ws.Range("K54").Value = 0
$ws.Range("L54").ClearContents()
$ws.Range("M54").Value = 8939.25
$ws.Range("N54").Value = -10057.25

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H55").Value = 7666.6665
$ws.Range("I55").Value = 3000
$ws.Range("J55").Value = 10000
$ws.Range("K55").Value = 9000
$ws.Range("L55").Value = 30000
$ws.Range("M55").Value = -8823
$ws.Range("N55").Value = -30354

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H75").Value = 0
$ws.Range("J75").Value = 0
$ws.Range("L75").ClearContents()
$ws.Range("N75").Value = 0

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H78").Value = 0
$ws.Range("J78").Value = 0
$ws.Range("L78").ClearContents()
$ws.Range("N78").Value = 0

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 793.7083
$ws.Range("I113").Value = 794.3889
$ws.Range("J113").Value = 791.6667
$ws.Range("K113").Value = 2383.1667
$ws.Range("L113").Value = 2375.0001
$ws.Range("M113").Value = -213.1667000000002
$ws.Range("N113").Value = -6715.0001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H135").Value = 8931305
$ws.Range("I135").Value = 553.6667
$ws.Range("J135").Value = 25006658
$ws.Range("K135").Value = 4983.0003
$ws.Range("L135").Value = 225059922
$ws.Range("M135").Value = -2448.0003
$ws.Range("N135").Value = -225064992

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 2125.5
$ws.Range("I113").Value = 2134.7368
$ws.Range("J113").Value = 2100.4285
$ws.Range("K113").Value = 2134.7368
$ws.Range("L113").Value = 2100.4285
$ws.Range("M113").Value = 35.26319999999987
$ws.Range("N113").Value = -6440.4285

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 3077.3333
$ws.Range("I126").Value = 1997.6666
$ws.Range("J126").Value = 4157
$ws.Range("K126").Value = 5992.9998
$ws.Range("L126").Value = 12471
$ws.Range("M126").Value = -3522.9998
$ws.Range("N126").Value = -17411

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H134").Value = 60163
$ws.Range("J134").Value = 60163
$ws.Range("L134").Value = 180489
$ws.Range("N134").Value = -185559

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H137").Value = 0
$ws.Range("J137").Value = 0
$ws.Range("L137").ClearContents()
$ws.Range("N137").Value = 0

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H114").Value = 26699
$ws.Range("J114").Value = 26699
$ws.Range("L114").Value = 26699
$ws.Range("N114").Value = -35377

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 3244.5
$ws.Range("I132").Value = 2923.2307
$ws.Range("J132").Value = 4079.8
$ws.Range("K132").Value = 8769.6921
$ws.Range("L132").Value = 12239.4
$ws.Range("M132").Value = -6239.6921
$ws.Range("N132").Value = -17299.4

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 4534.78
$ws.Range("I136").Value = 3028.2163
$ws.Range("J136").Value = 8822.691999999999
$ws.Range("K136").Value = 9084.6489
$ws.Range("L136").Value = 26468.076
$ws.Range("M136").Value = -6534.6489
$ws.Range("N136").Value = -31568.076

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H137").Value = 57800
$ws.Range("J137").Value = 57800
$ws.Range("L137").Value = 57800
$ws.Range("N137").Value = -68000

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H138").Value = 40726.125
$ws.Range("J138").Value = 40726.125
$ws.Range("L138").Value = 40726.125
$ws.Range("N138").Value = -51006.125
